$d = $word.ActiveDocument

# Locate the paragraph that contains the PCTO "patto formativo" sentence
# (the one that used to reference {A_SEDE_LEGALE} / {A_SEDE_SVOLGIMENTO}).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*residente in*A_SEDE_LEGALE*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="645CB60F" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordml" wp14:textId="25CE092A"><w:pPr><w:pStyle w:val="Normal" /><w:spacing w:line="360" w:lineRule="auto" /><w:jc w:val="both" /><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr></w:pPr><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t>Il/a sottoscritto/a {S_NOME} nato/a</w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t>…{</w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t>S_NATOA</w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t>}….</w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t>il</w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t>…{</w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t>S_NATOIL</w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t>}….</w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t>residente in {S_RESIDENZA} frequentante la classe…{S_CLASSE} {S_SEZIONE}, in procinto di frequentare attività dei PCTO nel period</w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t>o dal</w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /><w:b w:val="1" /><w:bCs w:val="1" /></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t xml:space="preserve"> {P_INIZIO} al {P_FINE} presso la struttura ospitante {A_NOME}, con sed</w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t xml:space="preserve">e legale  </w:t></w:r><w:r w:rsidRPr="31B07BA6" w:rsidR="31B07BA6"><w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond" w:eastAsia="Garamond" w:cs="Garamond" /></w:rPr><w:t>in {A_SEDE}, presso la sede di stage {A_SEDE}</w:t></w:r></w:p>
'@
    [void]$target.Range.InsertXML($xml)
}
